$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume/1h change (E) columns for rows with new values.
# D-column values that are plain decimal numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source data's text type)
# instead of auto-converting them to numbers.

$ws.Range("D2").Value = "39.832.01"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "2.227.19"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'292.01"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'87.31"
$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("D7").Value = "'0.514"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").Value = "'0.469"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").Value = "'30.69"
$ws.Range("E10").Value = "  +1.51%  "

$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("D12").Value = "'49.91"
$ws.Range("E12").Value = "  +5.25%  "

$ws.Range("E13").Value = "  +2.69%  "

$ws.Range("D14").Value = "'6.48"
$ws.Range("E14").Value = "  +2.62%  "

$ws.Range("D15").Value = "2.562.03"
$ws.Range("E15").Value = "  +0.73%  "

$ws.Range("D16").Value = "'13.86"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").Value = "2.217.15"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "'0.734"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "39.791.09"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  +1.21%  "

$ws.Range("D21").Value = "'11.16"
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").Value = "'5.76"
$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").Value = "'65.81"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").Value = "'237.49"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("D28").Value = "'23.19"
$ws.Range("E28").Value = "  +2.46%  "

$ws.Range("D29").Value = "'2.15"
$ws.Range("E29").Value = "  -2.11%  "

$ws.Range("D30").Value = "'9.27"
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").Value = "'156.75"
$ws.Range("E31").Value = "  +3.31%  "

$ws.Range("D32").Value = "'32.14"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").Value = "'4.99"
$ws.Range("E34").Value = "  +1.62%  "

$ws.Range("D35").Value = "'2.99"
$ws.Range("E35").Value = "  +7.70%  "

$ws.Range("D36").Value = "'0.0716"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("D38").Value = "'0.112"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D41").Value = "'15.41"
$ws.Range("E41").Value = "  -3.30%  "

$ws.Range("D42").Value = "2.114.50"
$ws.Range("E42").Value = "  +2.50%  "

$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "'18.44"
$ws.Range("E44").Value = "  +4.12%  "

$ws.Range("D45").Value = "'0.0272"
$ws.Range("E45").Value = "  +1.98%  "

$ws.Range("D46").Value = "'9.96"
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").Value = "'1.98"
$ws.Range("E47").Value = "  -8.77%  "

$ws.Range("D48").Value = "'2.73"
$ws.Range("E48").Value = "  +4.99%  "

$ws.Range("D49").Value = "2.434.20"
$ws.Range("E49").Value = "  +0.63%  "

$ws.Range("D50").Value = "'1.47"
$ws.Range("E50").Value = "  +3.57%  "

$ws.Range("E51").Value = "  +2.81%  "

# Rows 39 and 40: "Kaspa" and "ARBITRUM" swap places, each with refreshed
# price/volume figures.
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'1.76"
$ws.Range("E39").Value = "  +4.38%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.0991"
$ws.Range("E40").Value = "  +0.56%  "
